$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38 - shifts existing rows 38:90 down to 39:91
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with a new weekly record (same
# dimensions/origin/unit as the former row 38, new date + volume).
$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44495
$ws.Cells.Item(38, 5).Value = 4
$ws.Cells.Item(38, 6).Value = 100112044
$ws.Cells.Item(38, 7).Value = "Perejil"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 2860
$ws.Cells.Item(38, 11).Value = 1300
$ws.Cells.Item(38, 12).Value = 1500
$ws.Cells.Item(38, 13).Value = 1400
$ws.Cells.Item(38, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(38, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(38, 16).Value = 933
$ws.Cells.Item(38, 17).Value = 1.5
$ws.Cells.Item(38, 18).Value = "Hortaliza"
